$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "GETME"
$ws.Range("A5").Value = "ALSOGETME"
